$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: For every row that currently has a Notes cell in column B, move that
# content (and its style) into column C, then overwrite column B with the new CabNo
# header/value. Rows that are pure section-title rows (no existing B cell) are left
# untouched so they do not grow a phantom column-B cell.

# Row 3
$ws.Range("C3").Value = 'Notes'
$ws.Range("C3").HorizontalAlignment = -4147
$ws.Range("C3").IndentLevel = 2
$ws.Range("B3").Value = 'CabNo'
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").IndentLevel = 0

# Row 4
$ws.Range("C4").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C4").VerticalAlignment = -4160
$ws.Range("C4").WrapText = $true
$ws.Range("C4").IndentLevel = 2
$ws.Range("B4").Value = 'R1N1'
$ws.Range("B4").HorizontalAlignment = -4108
$ws.Range("B4").IndentLevel = 0
$ws.Range("B4").VerticalAlignment = -4160
$ws.Range("B4").WrapText = $true

# Row 5
$ws.Range("C5").Value = 'leave long, even if product width says otherwise'
$ws.Range("C5").VerticalAlignment = -4160
$ws.Range("C5").WrapText = $true
$ws.Range("C5").IndentLevel = 2
$ws.Range("B5").Value = 'R1N5'
$ws.Range("B5").HorizontalAlignment = -4108
$ws.Range("B5").IndentLevel = 0
$ws.Range("B5").VerticalAlignment = -4160
$ws.Range("B5").WrapText = $true

# Row 6
$ws.Range("C6").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C6").VerticalAlignment = -4160
$ws.Range("C6").WrapText = $true
$ws.Range("C6").IndentLevel = 2
$ws.Range("B6").Value = 'R1N2'
$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("B6").IndentLevel = 0
$ws.Range("B6").VerticalAlignment = -4160
$ws.Range("B6").WrapText = $true

# Row 7
$ws.Range("C7").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C7").VerticalAlignment = -4160
$ws.Range("C7").WrapText = $true
$ws.Range("C7").IndentLevel = 2
$ws.Range("B7").Value = 'R1N4'
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").IndentLevel = 0
$ws.Range("B7").VerticalAlignment = -4160
$ws.Range("B7").WrapText = $true

# Row 8
$ws.Range("C8").Value = 'leave long, even if product width says otherwise'
$ws.Range("C8").VerticalAlignment = -4160
$ws.Range("C8").WrapText = $true
$ws.Range("C8").IndentLevel = 2
$ws.Range("B8").Value = 'R1N9'
$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("B8").IndentLevel = 0
$ws.Range("B8").VerticalAlignment = -4160
$ws.Range("B8").WrapText = $true

# Row 11
$ws.Range("C11").Value = 'Notes'
$ws.Range("C11").HorizontalAlignment = -4147
$ws.Range("C11").IndentLevel = 2
$ws.Range("B11").Value = 'CabNo'
$ws.Range("B11").HorizontalAlignment = -4108
$ws.Range("B11").IndentLevel = 0

# Row 12
$ws.Range("C12").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C12").VerticalAlignment = -4160
$ws.Range("C12").WrapText = $true
$ws.Range("C12").IndentLevel = 2
$ws.Range("B12").Value = 'R2N1'
$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("B12").IndentLevel = 0
$ws.Range("B12").VerticalAlignment = -4160
$ws.Range("B12").WrapText = $true

# Row 13
$ws.Range("C13").Value = 'leave long, even if product width says otherwise'
$ws.Range("C13").VerticalAlignment = -4160
$ws.Range("C13").WrapText = $true
$ws.Range("C13").IndentLevel = 2
$ws.Range("B13").Value = 'R2N3'
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("B13").IndentLevel = 0
$ws.Range("B13").VerticalAlignment = -4160
$ws.Range("B13").WrapText = $true

# Row 16
$ws.Range("C16").Value = 'Notes'
$ws.Range("C16").HorizontalAlignment = -4147
$ws.Range("C16").IndentLevel = 2
$ws.Range("B16").Value = 'CabNo'
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("B16").IndentLevel = 0

# Row 17
$ws.Range("C17").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C17").VerticalAlignment = -4160
$ws.Range("C17").WrapText = $true
$ws.Range("C17").IndentLevel = 2
$ws.Range("B17").Value = 'R3N2'
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("B17").IndentLevel = 0
$ws.Range("B17").VerticalAlignment = -4160
$ws.Range("B17").WrapText = $true

# Row 18
$ws.Range("C18").Value = 'leave long, even if product width says otherwise'
$ws.Range("C18").VerticalAlignment = -4160
$ws.Range("C18").WrapText = $true
$ws.Range("C18").IndentLevel = 2
$ws.Range("B18").Value = 'R3N3'
$ws.Range("B18").HorizontalAlignment = -4108
$ws.Range("B18").IndentLevel = 0
$ws.Range("B18").VerticalAlignment = -4160
$ws.Range("B18").WrapText = $true

# Row 21
$ws.Range("C21").Value = 'Notes'
$ws.Range("C21").HorizontalAlignment = -4147
$ws.Range("C21").IndentLevel = 2
$ws.Range("B21").Value = 'CabNo'
$ws.Range("B21").HorizontalAlignment = -4108
$ws.Range("B21").IndentLevel = 0

# Row 22
$ws.Range("C22").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C22").VerticalAlignment = -4160
$ws.Range("C22").WrapText = $true
$ws.Range("C22").IndentLevel = 2
$ws.Range("B22").Value = 'R4N1'
$ws.Range("B22").HorizontalAlignment = -4108
$ws.Range("B22").IndentLevel = 0
$ws.Range("B22").VerticalAlignment = -4160
$ws.Range("B22").WrapText = $true

# Row 23
$ws.Range("C23").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C23").VerticalAlignment = -4160
$ws.Range("C23").WrapText = $true
$ws.Range("C23").IndentLevel = 2
$ws.Range("B23").Value = 'R4N4'
$ws.Range("B23").HorizontalAlignment = -4108
$ws.Range("B23").IndentLevel = 0
$ws.Range("B23").VerticalAlignment = -4160
$ws.Range("B23").WrapText = $true

# Row 24
$ws.Range("C24").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C24").VerticalAlignment = -4160
$ws.Range("C24").WrapText = $true
$ws.Range("C24").IndentLevel = 2
$ws.Range("B24").Value = 'R4N2'
$ws.Range("B24").HorizontalAlignment = -4108
$ws.Range("B24").IndentLevel = 0
$ws.Range("B24").VerticalAlignment = -4160
$ws.Range("B24").WrapText = $true

# Row 25
$ws.Range("C25").Value = 'leave long, even if product width says otherwise'
$ws.Range("C25").VerticalAlignment = -4160
$ws.Range("C25").WrapText = $true
$ws.Range("C25").IndentLevel = 2
$ws.Range("B25").Value = 'R4N5'
$ws.Range("B25").HorizontalAlignment = -4108
$ws.Range("B25").IndentLevel = 0
$ws.Range("B25").VerticalAlignment = -4160
$ws.Range("B25").WrapText = $true

# Row 26
$ws.Range("C26").Value = 'leave long, even if product width says otherwise'
$ws.Range("C26").VerticalAlignment = -4160
$ws.Range("C26").WrapText = $true
$ws.Range("C26").IndentLevel = 2
$ws.Range("B26").Value = 'R4N10'
$ws.Range("B26").HorizontalAlignment = -4108
$ws.Range("B26").IndentLevel = 0
$ws.Range("B26").VerticalAlignment = -4160
$ws.Range("B26").WrapText = $true

# Row 27
$ws.Range("C27").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C27").VerticalAlignment = -4160
$ws.Range("C27").WrapText = $true
$ws.Range("C27").IndentLevel = 2
$ws.Range("B27").Value = 'R4N8'
$ws.Range("B27").HorizontalAlignment = -4108
$ws.Range("B27").IndentLevel = 0
$ws.Range("B27").VerticalAlignment = -4160
$ws.Range("B27").WrapText = $true

# Row 28
$ws.Range("C28").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C28").VerticalAlignment = -4160
$ws.Range("C28").WrapText = $true
$ws.Range("C28").IndentLevel = 2
$ws.Range("B28").Value = 'R4N9'
$ws.Range("B28").HorizontalAlignment = -4108
$ws.Range("B28").IndentLevel = 0
$ws.Range("B28").VerticalAlignment = -4160
$ws.Range("B28").WrapText = $true

# Row 31
$ws.Range("C31").Value = 'Notes'
$ws.Range("C31").HorizontalAlignment = -4147
$ws.Range("C31").IndentLevel = 2
$ws.Range("B31").Value = 'CabNo'
$ws.Range("B31").HorizontalAlignment = -4108
$ws.Range("B31").IndentLevel = 0

# Row 32
$ws.Range("C32").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C32").VerticalAlignment = -4160
$ws.Range("C32").WrapText = $true
$ws.Range("C32").IndentLevel = 2
$ws.Range("B32").Value = 'R5N1'
$ws.Range("B32").HorizontalAlignment = -4108
$ws.Range("B32").IndentLevel = 0
$ws.Range("B32").VerticalAlignment = -4160
$ws.Range("B32").WrapText = $true

# Row 33
$ws.Range("C33").Value = 'leave long, even if product width says otherwise'
$ws.Range("C33").VerticalAlignment = -4160
$ws.Range("C33").WrapText = $true
$ws.Range("C33").IndentLevel = 2
$ws.Range("B33").Value = 'R5N3'
$ws.Range("B33").HorizontalAlignment = -4108
$ws.Range("B33").IndentLevel = 0
$ws.Range("B33").VerticalAlignment = -4160
$ws.Range("B33").WrapText = $true

# Row 36
$ws.Range("C36").Value = 'Notes'
$ws.Range("C36").HorizontalAlignment = -4147
$ws.Range("C36").IndentLevel = 2
$ws.Range("B36").Value = 'CabNo'
$ws.Range("B36").HorizontalAlignment = -4108
$ws.Range("B36").IndentLevel = 0

# Row 37
$ws.Range("C37").Value = 'leave long, even if product width says otherwise'
$ws.Range("C37").VerticalAlignment = -4160
$ws.Range("C37").WrapText = $true
$ws.Range("C37").IndentLevel = 2
$ws.Range("B37").Value = 'R6N3'
$ws.Range("B37").HorizontalAlignment = -4108
$ws.Range("B37").IndentLevel = 0
$ws.Range("B37").VerticalAlignment = -4160
$ws.Range("B37").WrapText = $true

# Row 40
$ws.Range("C40").Value = 'Notes'
$ws.Range("C40").HorizontalAlignment = -4147
$ws.Range("C40").IndentLevel = 2
$ws.Range("B40").Value = 'CabNo'
$ws.Range("B40").HorizontalAlignment = -4108
$ws.Range("B40").IndentLevel = 0

# Row 41
$ws.Range("C41").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C41").VerticalAlignment = -4160
$ws.Range("C41").WrapText = $true
$ws.Range("C41").IndentLevel = 2
$ws.Range("B41").Value = 'R7N1'
$ws.Range("B41").HorizontalAlignment = -4108
$ws.Range("B41").IndentLevel = 0
$ws.Range("B41").VerticalAlignment = -4160
$ws.Range("B41").WrapText = $true

# Row 42
$ws.Range("C42").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C42").VerticalAlignment = -4160
$ws.Range("C42").WrapText = $true
$ws.Range("C42").IndentLevel = 2
$ws.Range("B42").Value = 'R7N4'
$ws.Range("B42").HorizontalAlignment = -4108
$ws.Range("B42").IndentLevel = 0
$ws.Range("B42").VerticalAlignment = -4160
$ws.Range("B42").WrapText = $true

# Row 43
$ws.Range("C43").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C43").VerticalAlignment = -4160
$ws.Range("C43").WrapText = $true
$ws.Range("C43").IndentLevel = 2
$ws.Range("B43").Value = 'R7N2'
$ws.Range("B43").HorizontalAlignment = -4108
$ws.Range("B43").IndentLevel = 0
$ws.Range("B43").VerticalAlignment = -4160
$ws.Range("B43").WrapText = $true

# Row 44
$ws.Range("C44").Value = 'leave long, even if product width says otherwise'
$ws.Range("C44").VerticalAlignment = -4160
$ws.Range("C44").WrapText = $true
$ws.Range("C44").IndentLevel = 2
$ws.Range("B44").Value = 'R7N5'
$ws.Range("B44").HorizontalAlignment = -4108
$ws.Range("B44").IndentLevel = 0
$ws.Range("B44").VerticalAlignment = -4160
$ws.Range("B44").WrapText = $true

# Row 45
$ws.Range("C45").Value = 'leave long, even if product width says otherwise'
$ws.Range("C45").VerticalAlignment = -4160
$ws.Range("C45").WrapText = $true
$ws.Range("C45").IndentLevel = 2
$ws.Range("B45").Value = 'R7N10'
$ws.Range("B45").HorizontalAlignment = -4108
$ws.Range("B45").IndentLevel = 0
$ws.Range("B45").VerticalAlignment = -4160
$ws.Range("B45").WrapText = $true

# Row 46
$ws.Range("C46").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C46").VerticalAlignment = -4160
$ws.Range("C46").WrapText = $true
$ws.Range("C46").IndentLevel = 2
$ws.Range("B46").Value = 'R7N8'
$ws.Range("B46").HorizontalAlignment = -4108
$ws.Range("B46").IndentLevel = 0
$ws.Range("B46").VerticalAlignment = -4160
$ws.Range("B46").WrapText = $true

# Row 47
$ws.Range("C47").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C47").VerticalAlignment = -4160
$ws.Range("C47").WrapText = $true
$ws.Range("C47").IndentLevel = 2
$ws.Range("B47").Value = 'R7N9'
$ws.Range("B47").HorizontalAlignment = -4108
$ws.Range("B47").IndentLevel = 0
$ws.Range("B47").VerticalAlignment = -4160
$ws.Range("B47").WrapText = $true

# Row 50
$ws.Range("C50").Value = 'Notes'
$ws.Range("C50").HorizontalAlignment = -4147
$ws.Range("C50").IndentLevel = 2
$ws.Range("B50").Value = 'CabNo'
$ws.Range("B50").HorizontalAlignment = -4108
$ws.Range("B50").IndentLevel = 0

# Row 51
$ws.Range("C51").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C51").VerticalAlignment = -4160
$ws.Range("C51").WrapText = $true
$ws.Range("C51").IndentLevel = 2
$ws.Range("B51").Value = 'R8N1'
$ws.Range("B51").HorizontalAlignment = -4108
$ws.Range("B51").IndentLevel = 0
$ws.Range("B51").VerticalAlignment = -4160
$ws.Range("B51").WrapText = $true

# Row 52
$ws.Range("C52").Value = 'leave long, even if product width says otherwise'
$ws.Range("C52").VerticalAlignment = -4160
$ws.Range("C52").WrapText = $true
$ws.Range("C52").IndentLevel = 2
$ws.Range("B52").Value = 'R8N3'
$ws.Range("B52").HorizontalAlignment = -4108
$ws.Range("B52").IndentLevel = 0
$ws.Range("B52").VerticalAlignment = -4160
$ws.Range("B52").WrapText = $true

# Row 55
$ws.Range("C55").Value = 'Notes'
$ws.Range("C55").HorizontalAlignment = -4147
$ws.Range("C55").IndentLevel = 2
$ws.Range("B55").Value = 'CabNo'
$ws.Range("B55").HorizontalAlignment = -4108
$ws.Range("B55").IndentLevel = 0

# Row 56
$ws.Range("C56").Value = 'leave long, even if product width says otherwise'
$ws.Range("C56").VerticalAlignment = -4160
$ws.Range("C56").WrapText = $true
$ws.Range("C56").IndentLevel = 2
$ws.Range("B56").Value = 'R9C3'
$ws.Range("B56").HorizontalAlignment = -4108
$ws.Range("B56").IndentLevel = 0
$ws.Range("B56").VerticalAlignment = -4160
$ws.Range("B56").WrapText = $true

# Row 59
$ws.Range("C59").Value = 'Notes'
$ws.Range("C59").HorizontalAlignment = -4147
$ws.Range("C59").IndentLevel = 2
$ws.Range("B59").Value = 'CabNo'
$ws.Range("B59").HorizontalAlignment = -4108
$ws.Range("B59").IndentLevel = 0

# Row 60
$ws.Range("C60").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C60").VerticalAlignment = -4160
$ws.Range("C60").WrapText = $true
$ws.Range("C60").IndentLevel = 2
$ws.Range("B60").Value = 'R10N1'
$ws.Range("B60").HorizontalAlignment = -4108
$ws.Range("B60").IndentLevel = 0
$ws.Range("B60").VerticalAlignment = -4160
$ws.Range("B60").WrapText = $true

# Row 61
$ws.Range("C61").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C61").VerticalAlignment = -4160
$ws.Range("C61").WrapText = $true
$ws.Range("C61").IndentLevel = 2
$ws.Range("B61").Value = 'R10N4'
$ws.Range("B61").HorizontalAlignment = -4108
$ws.Range("B61").IndentLevel = 0
$ws.Range("B61").VerticalAlignment = -4160
$ws.Range("B61").WrapText = $true

# Row 62
$ws.Range("C62").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C62").VerticalAlignment = -4160
$ws.Range("C62").WrapText = $true
$ws.Range("C62").IndentLevel = 2
$ws.Range("B62").Value = 'R10N2'
$ws.Range("B62").HorizontalAlignment = -4108
$ws.Range("B62").IndentLevel = 0
$ws.Range("B62").VerticalAlignment = -4160
$ws.Range("B62").WrapText = $true

# Row 63
$ws.Range("C63").Value = 'leave long, even if product width says otherwise'
$ws.Range("C63").VerticalAlignment = -4160
$ws.Range("C63").WrapText = $true
$ws.Range("C63").IndentLevel = 2
$ws.Range("B63").Value = 'R10N5'
$ws.Range("B63").HorizontalAlignment = -4108
$ws.Range("B63").IndentLevel = 0
$ws.Range("B63").VerticalAlignment = -4160
$ws.Range("B63").WrapText = $true

# Row 64
$ws.Range("C64").Value = 'leave long, even if product width says otherwise'
$ws.Range("C64").VerticalAlignment = -4160
$ws.Range("C64").WrapText = $true
$ws.Range("C64").IndentLevel = 2
$ws.Range("B64").Value = 'R10N10'
$ws.Range("B64").HorizontalAlignment = -4108
$ws.Range("B64").IndentLevel = 0
$ws.Range("B64").VerticalAlignment = -4160
$ws.Range("B64").WrapText = $true

# Row 65
$ws.Range("C65").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C65").VerticalAlignment = -4160
$ws.Range("C65").WrapText = $true
$ws.Range("C65").IndentLevel = 2
$ws.Range("B65").Value = 'R10N8'
$ws.Range("B65").HorizontalAlignment = -4108
$ws.Range("B65").IndentLevel = 0
$ws.Range("B65").VerticalAlignment = -4160
$ws.Range("B65").WrapText = $true

# Row 66
$ws.Range("C66").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C66").VerticalAlignment = -4160
$ws.Range("C66").WrapText = $true
$ws.Range("C66").IndentLevel = 2
$ws.Range("B66").Value = 'R10N9'
$ws.Range("B66").HorizontalAlignment = -4108
$ws.Range("B66").IndentLevel = 0
$ws.Range("B66").VerticalAlignment = -4160
$ws.Range("B66").WrapText = $true

# Row 69
$ws.Range("C69").Value = 'Notes'
$ws.Range("C69").HorizontalAlignment = -4147
$ws.Range("C69").IndentLevel = 2
$ws.Range("B69").Value = 'CabNo'
$ws.Range("B69").HorizontalAlignment = -4108
$ws.Range("B69").IndentLevel = 0

# Row 70
$ws.Range("C70").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C70").VerticalAlignment = -4160
$ws.Range("C70").WrapText = $true
$ws.Range("C70").IndentLevel = 2
$ws.Range("B70").Value = 'R11N1'
$ws.Range("B70").HorizontalAlignment = -4108
$ws.Range("B70").IndentLevel = 0
$ws.Range("B70").VerticalAlignment = -4160
$ws.Range("B70").WrapText = $true

# Row 71
$ws.Range("C71").Value = 'leave long, even if product width says otherwise'
$ws.Range("C71").VerticalAlignment = -4160
$ws.Range("C71").WrapText = $true
$ws.Range("C71").IndentLevel = 2
$ws.Range("B71").Value = 'R11N3'
$ws.Range("B71").HorizontalAlignment = -4108
$ws.Range("B71").IndentLevel = 0
$ws.Range("B71").VerticalAlignment = -4160
$ws.Range("B71").WrapText = $true

# Row 74
$ws.Range("C74").Value = 'Notes'
$ws.Range("C74").HorizontalAlignment = -4147
$ws.Range("C74").IndentLevel = 2
$ws.Range("B74").Value = 'CabNo'
$ws.Range("B74").HorizontalAlignment = -4108
$ws.Range("B74").IndentLevel = 0

# Row 75
$ws.Range("C75").Value = 'leave long, even if product width says otherwise'
$ws.Range("C75").VerticalAlignment = -4160
$ws.Range("C75").WrapText = $true
$ws.Range("C75").IndentLevel = 2
$ws.Range("B75").Value = 'R12N3'
$ws.Range("B75").HorizontalAlignment = -4108
$ws.Range("B75").IndentLevel = 0
$ws.Range("B75").VerticalAlignment = -4160
$ws.Range("B75").WrapText = $true

# Row 78
$ws.Range("C78").Value = 'Notes'
$ws.Range("C78").HorizontalAlignment = -4147
$ws.Range("C78").IndentLevel = 2
$ws.Range("B78").Value = 'CabNo'
$ws.Range("B78").HorizontalAlignment = -4108
$ws.Range("B78").IndentLevel = 0

# Row 79
$ws.Range("C79").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C79").VerticalAlignment = -4160
$ws.Range("C79").WrapText = $true
$ws.Range("C79").IndentLevel = 2
$ws.Range("B79").Value = 'R13N1'
$ws.Range("B79").HorizontalAlignment = -4108
$ws.Range("B79").IndentLevel = 0
$ws.Range("B79").VerticalAlignment = -4160
$ws.Range("B79").WrapText = $true

# Row 80
$ws.Range("C80").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C80").VerticalAlignment = -4160
$ws.Range("C80").WrapText = $true
$ws.Range("C80").IndentLevel = 2
$ws.Range("B80").Value = 'R13N4'
$ws.Range("B80").HorizontalAlignment = -4108
$ws.Range("B80").IndentLevel = 0
$ws.Range("B80").VerticalAlignment = -4160
$ws.Range("B80").WrapText = $true

# Row 81
$ws.Range("C81").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C81").VerticalAlignment = -4160
$ws.Range("C81").WrapText = $true
$ws.Range("C81").IndentLevel = 2
$ws.Range("B81").Value = 'R13N2'
$ws.Range("B81").HorizontalAlignment = -4108
$ws.Range("B81").IndentLevel = 0
$ws.Range("B81").VerticalAlignment = -4160
$ws.Range("B81").WrapText = $true

# Row 82
$ws.Range("C82").Value = 'leave long, even if product width says otherwise'
$ws.Range("C82").VerticalAlignment = -4160
$ws.Range("C82").WrapText = $true
$ws.Range("C82").IndentLevel = 2
$ws.Range("B82").Value = 'R13N5'
$ws.Range("B82").HorizontalAlignment = -4108
$ws.Range("B82").IndentLevel = 0
$ws.Range("B82").VerticalAlignment = -4160
$ws.Range("B82").WrapText = $true

# Row 83
$ws.Range("C83").Value = 'leave long, even if product width says otherwise'
$ws.Range("C83").VerticalAlignment = -4160
$ws.Range("C83").WrapText = $true
$ws.Range("C83").IndentLevel = 2
$ws.Range("B83").Value = 'R13N10'
$ws.Range("B83").HorizontalAlignment = -4108
$ws.Range("B83").IndentLevel = 0
$ws.Range("B83").VerticalAlignment = -4160
$ws.Range("B83").WrapText = $true

# Row 84
$ws.Range("C84").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C84").VerticalAlignment = -4160
$ws.Range("C84").WrapText = $true
$ws.Range("C84").IndentLevel = 2
$ws.Range("B84").Value = 'R13N8'
$ws.Range("B84").HorizontalAlignment = -4108
$ws.Range("B84").IndentLevel = 0
$ws.Range("B84").VerticalAlignment = -4160
$ws.Range("B84").WrapText = $true

# Row 85
$ws.Range("C85").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C85").VerticalAlignment = -4160
$ws.Range("C85").WrapText = $true
$ws.Range("C85").IndentLevel = 2
$ws.Range("B85").Value = 'R13N9'
$ws.Range("B85").HorizontalAlignment = -4108
$ws.Range("B85").IndentLevel = 0
$ws.Range("B85").VerticalAlignment = -4160
$ws.Range("B85").WrapText = $true

# Row 88
$ws.Range("C88").Value = 'Notes'
$ws.Range("C88").HorizontalAlignment = -4147
$ws.Range("C88").IndentLevel = 2
$ws.Range("B88").Value = 'CabNo'
$ws.Range("B88").HorizontalAlignment = -4108
$ws.Range("B88").IndentLevel = 0

# Row 89
$ws.Range("C89").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C89").VerticalAlignment = -4160
$ws.Range("C89").WrapText = $true
$ws.Range("C89").IndentLevel = 2
$ws.Range("B89").Value = 'R14N1'
$ws.Range("B89").HorizontalAlignment = -4108
$ws.Range("B89").IndentLevel = 0
$ws.Range("B89").VerticalAlignment = -4160
$ws.Range("B89").WrapText = $true

# Row 90
$ws.Range("C90").Value = 'leave long, even if product width says otherwise'
$ws.Range("C90").VerticalAlignment = -4160
$ws.Range("C90").WrapText = $true
$ws.Range("C90").IndentLevel = 2
$ws.Range("B90").Value = 'R14N3'
$ws.Range("B90").HorizontalAlignment = -4108
$ws.Range("B90").IndentLevel = 0
$ws.Range("B90").VerticalAlignment = -4160
$ws.Range("B90").WrapText = $true

# Row 93
$ws.Range("C93").Value = 'Notes'
$ws.Range("C93").HorizontalAlignment = -4147
$ws.Range("C93").IndentLevel = 2
$ws.Range("B93").Value = 'CabNo'
$ws.Range("B93").HorizontalAlignment = -4108
$ws.Range("B93").IndentLevel = 0

# Row 94
$ws.Range("C94").Value = 'leave long, even if product width says otherwise'
$ws.Range("C94").VerticalAlignment = -4160
$ws.Range("C94").WrapText = $true
$ws.Range("C94").IndentLevel = 2
$ws.Range("B94").Value = 'R15N3'
$ws.Range("B94").HorizontalAlignment = -4108
$ws.Range("B94").IndentLevel = 0
$ws.Range("B94").VerticalAlignment = -4160
$ws.Range("B94").WrapText = $true

# Row 97
$ws.Range("C97").Value = 'Notes'
$ws.Range("C97").HorizontalAlignment = -4147
$ws.Range("C97").IndentLevel = 2
$ws.Range("B97").Value = 'CabNo'
$ws.Range("B97").HorizontalAlignment = -4108
$ws.Range("B97").IndentLevel = 0

# Row 98
$ws.Range("C98").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C98").VerticalAlignment = -4160
$ws.Range("C98").WrapText = $true
$ws.Range("C98").IndentLevel = 2
$ws.Range("B98").Value = 'R16N1'
$ws.Range("B98").HorizontalAlignment = -4108
$ws.Range("B98").IndentLevel = 0
$ws.Range("B98").VerticalAlignment = -4160
$ws.Range("B98").WrapText = $true

# Row 99
$ws.Range("C99").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C99").VerticalAlignment = -4160
$ws.Range("C99").WrapText = $true
$ws.Range("C99").IndentLevel = 2
$ws.Range("B99").Value = 'R16N4'
$ws.Range("B99").HorizontalAlignment = -4108
$ws.Range("B99").IndentLevel = 0
$ws.Range("B99").VerticalAlignment = -4160
$ws.Range("B99").WrapText = $true

# Row 100
$ws.Range("C100").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C100").VerticalAlignment = -4160
$ws.Range("C100").WrapText = $true
$ws.Range("C100").IndentLevel = 2
$ws.Range("B100").Value = 'R16N2'
$ws.Range("B100").HorizontalAlignment = -4108
$ws.Range("B100").IndentLevel = 0
$ws.Range("B100").VerticalAlignment = -4160
$ws.Range("B100").WrapText = $true

# Row 101
$ws.Range("C101").Value = 'leave long, even if product width says otherwise'
$ws.Range("C101").VerticalAlignment = -4160
$ws.Range("C101").WrapText = $true
$ws.Range("C101").IndentLevel = 2
$ws.Range("B101").Value = 'R16N5'
$ws.Range("B101").HorizontalAlignment = -4108
$ws.Range("B101").IndentLevel = 0
$ws.Range("B101").VerticalAlignment = -4160
$ws.Range("B101").WrapText = $true

# Row 102
$ws.Range("C102").Value = 'leave long, even if product width says otherwise'
$ws.Range("C102").VerticalAlignment = -4160
$ws.Range("C102").WrapText = $true
$ws.Range("C102").IndentLevel = 2
$ws.Range("B102").Value = 'R16N10'
$ws.Range("B102").HorizontalAlignment = -4108
$ws.Range("B102").IndentLevel = 0
$ws.Range("B102").VerticalAlignment = -4160
$ws.Range("B102").WrapText = $true

# Row 103
$ws.Range("C103").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C103").VerticalAlignment = -4160
$ws.Range("C103").WrapText = $true
$ws.Range("C103").IndentLevel = 2
$ws.Range("B103").Value = 'R16N8'
$ws.Range("B103").HorizontalAlignment = -4108
$ws.Range("B103").IndentLevel = 0
$ws.Range("B103").VerticalAlignment = -4160
$ws.Range("B103").WrapText = $true

# Row 104
$ws.Range("C104").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C104").VerticalAlignment = -4160
$ws.Range("C104").WrapText = $true
$ws.Range("C104").IndentLevel = 2
$ws.Range("B104").Value = 'R16N9'
$ws.Range("B104").HorizontalAlignment = -4108
$ws.Range("B104").IndentLevel = 0
$ws.Range("B104").VerticalAlignment = -4160
$ws.Range("B104").WrapText = $true

# Row 107
$ws.Range("C107").Value = 'Notes'
$ws.Range("C107").HorizontalAlignment = -4147
$ws.Range("C107").IndentLevel = 2
$ws.Range("B107").Value = 'CabNo'
$ws.Range("B107").HorizontalAlignment = -4108
$ws.Range("B107").IndentLevel = 0

# Row 108
$ws.Range("C108").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C108").VerticalAlignment = -4160
$ws.Range("C108").WrapText = $true
$ws.Range("C108").IndentLevel = 2
$ws.Range("B108").Value = 'R17N1'
$ws.Range("B108").HorizontalAlignment = -4108
$ws.Range("B108").IndentLevel = 0
$ws.Range("B108").VerticalAlignment = -4160
$ws.Range("B108").WrapText = $true

# Row 109
$ws.Range("C109").Value = 'leave long, even if product width says otherwise'
$ws.Range("C109").VerticalAlignment = -4160
$ws.Range("C109").WrapText = $true
$ws.Range("C109").IndentLevel = 2
$ws.Range("B109").Value = 'R17N3'
$ws.Range("B109").HorizontalAlignment = -4108
$ws.Range("B109").IndentLevel = 0
$ws.Range("B109").VerticalAlignment = -4160
$ws.Range("B109").WrapText = $true

# Row 112
$ws.Range("C112").Value = 'Notes'
$ws.Range("C112").HorizontalAlignment = -4147
$ws.Range("C112").IndentLevel = 2
$ws.Range("B112").Value = 'CabNo'
$ws.Range("B112").HorizontalAlignment = -4108
$ws.Range("B112").IndentLevel = 0

# Row 113
$ws.Range("C113").Value = 'leave long, even if product width says otherwise'
$ws.Range("C113").VerticalAlignment = -4160
$ws.Range("C113").WrapText = $true
$ws.Range("C113").IndentLevel = 2
$ws.Range("B113").Value = 'R18N3'
$ws.Range("B113").HorizontalAlignment = -4108
$ws.Range("B113").IndentLevel = 0
$ws.Range("B113").VerticalAlignment = -4160
$ws.Range("B113").WrapText = $true

# Row 116
$ws.Range("C116").Value = 'Notes'
$ws.Range("C116").HorizontalAlignment = -4147
$ws.Range("C116").IndentLevel = 2
$ws.Range("B116").Value = 'CabNo'
$ws.Range("B116").HorizontalAlignment = -4108
$ws.Range("B116").IndentLevel = 0

# Row 117
$ws.Range("C117").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C117").VerticalAlignment = -4160
$ws.Range("C117").WrapText = $true
$ws.Range("C117").IndentLevel = 2
$ws.Range("B117").Value = 'R19N1'
$ws.Range("B117").HorizontalAlignment = -4108
$ws.Range("B117").IndentLevel = 0
$ws.Range("B117").VerticalAlignment = -4160
$ws.Range("B117").WrapText = $true

# Row 118
$ws.Range("C118").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C118").VerticalAlignment = -4160
$ws.Range("C118").WrapText = $true
$ws.Range("C118").IndentLevel = 2
$ws.Range("B118").Value = 'R19N4'
$ws.Range("B118").HorizontalAlignment = -4108
$ws.Range("B118").IndentLevel = 0
$ws.Range("B118").VerticalAlignment = -4160
$ws.Range("B118").WrapText = $true

# Row 119
$ws.Range("C119").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C119").VerticalAlignment = -4160
$ws.Range("C119").WrapText = $true
$ws.Range("C119").IndentLevel = 2
$ws.Range("B119").Value = 'R19N2'
$ws.Range("B119").HorizontalAlignment = -4108
$ws.Range("B119").IndentLevel = 0
$ws.Range("B119").VerticalAlignment = -4160
$ws.Range("B119").WrapText = $true

# Row 120
$ws.Range("C120").Value = 'leave long, even if product width says otherwise'
$ws.Range("C120").VerticalAlignment = -4160
$ws.Range("C120").WrapText = $true
$ws.Range("C120").IndentLevel = 2
$ws.Range("B120").Value = 'R19N5'
$ws.Range("B120").HorizontalAlignment = -4108
$ws.Range("B120").IndentLevel = 0
$ws.Range("B120").VerticalAlignment = -4160
$ws.Range("B120").WrapText = $true

# Row 121
$ws.Range("C121").Value = 'leave long, even if product width says otherwise'
$ws.Range("C121").VerticalAlignment = -4160
$ws.Range("C121").WrapText = $true
$ws.Range("C121").IndentLevel = 2
$ws.Range("B121").Value = 'R19N10'
$ws.Range("B121").HorizontalAlignment = -4108
$ws.Range("B121").IndentLevel = 0
$ws.Range("B121").VerticalAlignment = -4160
$ws.Range("B121").WrapText = $true

# Row 122
$ws.Range("C122").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C122").VerticalAlignment = -4160
$ws.Range("C122").WrapText = $true
$ws.Range("C122").IndentLevel = 2
$ws.Range("B122").Value = 'R19N8'
$ws.Range("B122").HorizontalAlignment = -4108
$ws.Range("B122").IndentLevel = 0
$ws.Range("B122").VerticalAlignment = -4160
$ws.Range("B122").WrapText = $true

# Row 123
$ws.Range("C123").Value = '"L" shaped product with face height equal to Height dimension.
Use cutoff for cleat. Pocketscrew cleat on B side.
Attach cleat to face in shop with pocket holes on top side.
Leave long, even if product width says otherwise.'
$ws.Range("C123").VerticalAlignment = -4160
$ws.Range("C123").WrapText = $true
$ws.Range("C123").IndentLevel = 2
$ws.Range("B123").Value = 'R19N9'
$ws.Range("B123").HorizontalAlignment = -4108
$ws.Range("B123").IndentLevel = 0
$ws.Range("B123").VerticalAlignment = -4160
$ws.Range("B123").WrapText = $true

# Row 126
$ws.Range("C126").Value = 'Notes'
$ws.Range("C126").HorizontalAlignment = -4147
$ws.Range("C126").IndentLevel = 2
$ws.Range("B126").Value = 'CabNo'
$ws.Range("B126").HorizontalAlignment = -4108
$ws.Range("B126").IndentLevel = 0

# Row 127
$ws.Range("C127").Value = 'Pocketscrew B-side of cleat every ~8"'
$ws.Range("C127").VerticalAlignment = -4160
$ws.Range("C127").WrapText = $true
$ws.Range("C127").IndentLevel = 2
$ws.Range("B127").Value = 'R20N1'
$ws.Range("B127").HorizontalAlignment = -4108
$ws.Range("B127").IndentLevel = 0
$ws.Range("B127").VerticalAlignment = -4160
$ws.Range("B127").WrapText = $true

# Row 128
$ws.Range("C128").Value = 'leave long, even if product width says otherwise'
$ws.Range("C128").VerticalAlignment = -4160
$ws.Range("C128").WrapText = $true
$ws.Range("C128").IndentLevel = 2
$ws.Range("B128").Value = 'R20N3'
$ws.Range("B128").HorizontalAlignment = -4108
$ws.Range("B128").IndentLevel = 0
$ws.Range("B128").VerticalAlignment = -4160
$ws.Range("B128").WrapText = $true

# Row 131
$ws.Range("C131").Value = 'Notes'
$ws.Range("C131").HorizontalAlignment = -4147
$ws.Range("C131").IndentLevel = 2
$ws.Range("B131").Value = 'CabNo'
$ws.Range("B131").HorizontalAlignment = -4108
$ws.Range("B131").IndentLevel = 0

# Row 132
$ws.Range("C132").Value = 'leave long, even if product width says otherwise'
$ws.Range("C132").VerticalAlignment = -4160
$ws.Range("C132").WrapText = $true
$ws.Range("C132").IndentLevel = 2
$ws.Range("B132").Value = 'R21N3'
$ws.Range("B132").HorizontalAlignment = -4108
$ws.Range("B132").IndentLevel = 0
$ws.Range("B132").VerticalAlignment = -4160
$ws.Range("B132").WrapText = $true

# --- Step 2: extend merged section-title cells from column B to column C
$ws.Range("A1:B1").UnMerge()
$ws.Range("A1:C1").Merge()
$ws.Range("A2:B2").UnMerge()
$ws.Range("A2:C2").Merge()
$ws.Range("A10:B10").UnMerge()
$ws.Range("A10:C10").Merge()
$ws.Range("A15:B15").UnMerge()
$ws.Range("A15:C15").Merge()
$ws.Range("A20:B20").UnMerge()
$ws.Range("A20:C20").Merge()
$ws.Range("A30:B30").UnMerge()
$ws.Range("A30:C30").Merge()
$ws.Range("A35:B35").UnMerge()
$ws.Range("A35:C35").Merge()
$ws.Range("A39:B39").UnMerge()
$ws.Range("A39:C39").Merge()
$ws.Range("A49:B49").UnMerge()
$ws.Range("A49:C49").Merge()
$ws.Range("A54:B54").UnMerge()
$ws.Range("A54:C54").Merge()
$ws.Range("A58:B58").UnMerge()
$ws.Range("A58:C58").Merge()
$ws.Range("A68:B68").UnMerge()
$ws.Range("A68:C68").Merge()
$ws.Range("A73:B73").UnMerge()
$ws.Range("A73:C73").Merge()
$ws.Range("A77:B77").UnMerge()
$ws.Range("A77:C77").Merge()
$ws.Range("A87:B87").UnMerge()
$ws.Range("A87:C87").Merge()
$ws.Range("A92:B92").UnMerge()
$ws.Range("A92:C92").Merge()
$ws.Range("A96:B96").UnMerge()
$ws.Range("A96:C96").Merge()
$ws.Range("A106:B106").UnMerge()
$ws.Range("A106:C106").Merge()
$ws.Range("A111:B111").UnMerge()
$ws.Range("A111:C111").Merge()
$ws.Range("A115:B115").UnMerge()
$ws.Range("A115:C115").Merge()
$ws.Range("A125:B125").UnMerge()
$ws.Range("A125:C125").Merge()
$ws.Range("A130:B130").UnMerge()
$ws.Range("A130:C130").Merge()

# --- Step 3: column widths.  Column A (width 25) stays as-is.  Column B becomes the
# narrow new CabNo column (target XML width = 8) and column C becomes the old Notes
# column (target XML width = 50).  This engine's ColumnWidth setter adds a constant
# 0.8333333 (5/6) padding versus the raw XML "width" attribute, so we subtract it back
# out to land on an exact integer width in the saved file.
$ws.Columns("B").ColumnWidth = 8 - 0.8333333
$ws.Columns("C").ColumnWidth = 50 - 0.8333333

# --- Step 4: Print area grows from A1:B134 to A1:C134
$ws.PageSetup.PrintArea = "`$A`$1:`$C`$134"
